$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.993.43'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.641.12'
$ws.Range('E3').Value = '  -0.45%  '
$cell = $ws.Range('D4')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = $origStyle
$ws.Range('E4').Value = '  -1.17%  '
$ws.Range('D5').Value = '214.89'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = '0.5083'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('D8').Value = '0.2585'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = '0.06362'
$ws.Range('E9').Value = '  -1.06%  '
$ws.Range('D10').Value = '19.87'
$ws.Range('E10').Value = '  +0.92%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.07750'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('D12').Value = '4.292'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '1.642.23'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = '0.5489'
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0₅7753'
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '64.40'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = '26.005.40'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = '196.72'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').Value = '4.455'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').Value = '9.969'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').Value = '6.121'
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('E23').Value = '  -0.69%  '
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.890'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').Value = '142.97'
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('D26').Value = '0.1258'
$ws.Range('E26').Value = '  +9.31%  '
$ws.Range('D27').Value = '6.877'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = '15.66'
$ws.Range('E28').Value = '  -0.68%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.240'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('D30').Value = '0.04907'
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('D31').Value = '3.278'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Value = '3.219'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = '1.552'
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').Value = '2.374'
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').Value = '0.9182'
$ws.Range('E35').Value = '  +2.51%  '
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.570'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  -1.12%  '
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.5550'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').Value = '1.100.41'
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('D39').Value = '0.01573'
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range('E40').Value = '  -0.74%  '
$ws.Range('D41').Value = '5.613'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').Value = '0.8037'
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('D43').Value = '98.86'
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('E44').Value = '  -4.02%  '
$ws.Range('D45').Value = '1.784.42'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = '0.4533'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').Value = '55.31'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').Value = '0.05193'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('D50').Value = '7.552'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').Value = '1.003'
$ws.Range('E51').Value = '  -0.48%  '
